$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 136.58333
$ws.Range("J28").Value = 300
$ws.Range("L28").Value = 300
$ws.Range("N28").Value = -1270

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1541.6875
$ws.Range("J112").Value = 1611.1333
$ws.Range("L112").Value = 4833.3999
$ws.Range("N112").Value = -7049.3999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 399.66666
$ws.Range("I118").Value = 399.66666
$ws.Range("K118").Value = 1198.99998
$ws.Range("M118").Value = 458.0000199999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1461.4348
$ws.Range("J137").Value = 1983.5
$ws.Range("L137").Value = 5950.5
$ws.Range("N137").Value = -11050.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1656.6666
$ws.Range("I110").Value = 1478.5555
$ws.Range("J110").Value = 2191
$ws.Range("K110").Value = 1478.5555
$ws.Range("L110").Value = 2191
$ws.Range("M110").Value = 566.4445000000001
$ws.Range("N110").Value = -6281

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 15501.25
$ws.Range("I132").Value = 26003
$ws.Range("K132").Value = 78009
$ws.Range("M132").Value = -75479

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -11058

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 28998.334
$ws.Range("J88").Value = 28998.334
$ws.Range("L88").Value = 28998.334
$ws.Range("N88").Value = -29810.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 28998.334
$ws.Range("J91").Value = 28998.334
$ws.Range("L91").Value = 28998.334
$ws.Range("N91").Value = -31806.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2500
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 79255.84
$ws.Range("I134").Value = 201462.4
$ws.Range("J134").Value = 2876.75
$ws.Range("K134").Value = 604387.2
$ws.Range("L134").Value = 8630.25
$ws.Range("M134").Value = -601852.2
$ws.Range("N134").Value = -13700.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1100.2
$ws.Range("I16").Value = 849.25
$ws.Range("J16").Value = 2104
$ws.Range("K16").Value = 849.25
$ws.Range("L16").Value = 2104
$ws.Range("M16").Value = -562.25
$ws.Range("N16").Value = -2678

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8336225.5
$ws.Range("I31").Value = 2973.182
$ws.Range("J31").Value = 100002000
$ws.Range("K31").Value = 2973.182
$ws.Range("L31").Value = 100002000
$ws.Range("M31").Value = -2678.182
$ws.Range("N31").Value = -100002590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8336225.5
$ws.Range("I34").Value = 2973.182
$ws.Range("J34").Value = 100002000
$ws.Range("K34").Value = 2973.182
$ws.Range("L34").Value = 100002000
$ws.Range("M34").Value = -2771.182
$ws.Range("N34").Value = -100002404

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 521.0454999999999
$ws.Range("I107").Value = 542.4666999999999
$ws.Range("J107").Value = 475.14285
$ws.Range("K107").Value = 542.4666999999999
$ws.Range("L107").Value = 475.14285
$ws.Range("M107").Value = 1377.5333
$ws.Range("N107").Value = -4315.14285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1100.2
$ws.Range("I113").Value = 849.25
$ws.Range("J113").Value = 2104
$ws.Range("K113").Value = 849.25
$ws.Range("L113").Value = 2104
$ws.Range("M113").Value = 1320.75
$ws.Range("N113").Value = -6444

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2950.3333
$ws.Range("I132").Value = 2784.2666
$ws.Range("K132").Value = 8352.799800000001
$ws.Range("M132").Value = -5822.799800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 742701.4399999999
$ws.Range("I122").Value = 2567.1667
$ws.Range("J122").Value = 1236124.4
$ws.Range("K122").Value = 23104.5003
$ws.Range("L122").Value = 11125119.6
$ws.Range("M122").Value = -20654.5003
$ws.Range("N122").Value = -11130019.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3401337
$ws.Range("I131").Value = 9575.833000000001
$ws.Range("J131").Value = 5339486.5
$ws.Range("K131").Value = 28727.499
$ws.Range("L131").Value = 16018459.5
$ws.Range("M131").Value = -23687.499
$ws.Range("N131").Value = -16028539.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1138.1212
$ws.Range("I132").Value = 484.15
$ws.Range("K132").Value = 4357.349999999999
$ws.Range("M132").Value = -1827.349999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 66151596
$ws.Range("I137").Value = 47619748
$ws.Range("J137").Value = 75417520
$ws.Range("K137").Value = 142859244
$ws.Range("L137").Value = 226252560
$ws.Range("M137").Value = -142854144
$ws.Range("N137").Value = -226262760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 537.9
$ws.Range("I97").Value = 369.73334
$ws.Range("J97").Value = 706.06665
$ws.Range("K97").Value = 369.73334
$ws.Range("L97").Value = 706.06665
$ws.Range("M97").Value = 126.26666
$ws.Range("N97").Value = -1698.06665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 675.4
$ws.Range("J107").Value = 908.6
$ws.Range("L107").Value = 908.6
$ws.Range("N107").Value = -4748.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 202613.7
$ws.Range("I132").Value = 335454
$ws.Range("J132").Value = 3353.25
$ws.Range("K132").Value = 1006362
$ws.Range("L132").Value = 10059.75
$ws.Range("M132").Value = -1003832
$ws.Range("N132").Value = -15119.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 8333.333000000001
$ws.Range("I20").Value = 5000
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = -4774
$ws.Range("N20").Value = -10452

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 13603.311
$ws.Range("I132").Value = 29976.727
$ws.Range("J132").Value = 3597.3333
$ws.Range("K132").Value = 89930.181
$ws.Range("L132").Value = 10791.9999
$ws.Range("M132").Value = -87400.181
$ws.Range("N132").Value = -15851.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6740.2
$ws.Range("I136").Value = 21040.8
$ws.Range("J136").Value = 1973.3334
$ws.Range("K136").Value = 63122.39999999999
$ws.Range("L136").Value = 5920.0002
$ws.Range("M136").Value = -60572.39999999999
$ws.Range("N136").Value = -11020.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1056.35
$ws.Range("I122").Value = 1042.5555
$ws.Range("J122").Value = 1067.6364
$ws.Range("K122").Value = 3127.6665
$ws.Range("L122").Value = 3202.9092
$ws.Range("M122").Value = -677.6664999999998
$ws.Range("N122").Value = -8102.9092
